# Applies the cryptos-list refresh described by the commit diff:
# updates Price (D) / Volume-1h (E) values for most rows, and swaps
# the Aave / InjectiveProtocol rows (42-43) including their B/C text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.895.25'
$ws.Range("E2").Value = '  -0.59%  '
$ws.Range("D3").Value = '2.041.11'
$ws.Range("E3").Value = '  -0.95%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = "'" + '227.93'
$ws.Range("E5").Value = '  -1.10%  '
$ws.Range("D6").Value = "'" + '0.615'
$ws.Range("E6").Value = '  -0.19%  '
$ws.Range("D7").Value = "'" + '60.46'
$ws.Range("E7").Value = '  +4.02%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("E9").Value = '  +0.22%  '
$ws.Range("D10").Value = "'" + '0.0816'
$ws.Range("E10").Value = '  +1.27%  '
$ws.Range("E11").Value = '  +0.57%  '
$ws.Range("D12").Value = "'" + '14.70'
$ws.Range("E12").Value = '  +0.48%  '
$ws.Range("D13").Value = '2.344.28'
$ws.Range("E13").Value = '  -0.89%  '
$ws.Range("D14").Value = "'" + '21.10'
$ws.Range("E14").Value = '  +2.09%  '
$ws.Range("D15").Value = "'" + '0.765'
$ws.Range("E15").Value = '  +1.42%  '
$ws.Range("D16").Value = "'" + '5.20'
$ws.Range("E16").Value = '  -1.91%  '
$ws.Range("D17").Value = '2.033.92'
$ws.Range("E17").Value = '  -0.92%  '
$ws.Range("D18").Value = '37.856.26'
$ws.Range("E18").Value = '  -0.46%  '
$ws.Range("E19").Value = '  -1.49%  '
$ws.Range("D20").Value = "'" + '69.90'
$ws.Range("E20").Value = '  +0.30%  '
$ws.Range("E21").Value = '  -0.64%  '
$ws.Range("D22").Value = "'" + '225.49'
$ws.Range("E22").Value = '  +0.45%  '
$ws.Range("E23").Value = '  -0.05%  '
$ws.Range("E24").Value = '  -2.44%  '
$ws.Range("D25").Value = "'" + '2.22'
$ws.Range("E25").Value = '  -1.78%  '
$ws.Range("D26").Value = "'" + '9.26'
$ws.Range("E26").Value = '  -0.68%  '
$ws.Range("D27").Value = "'" + '165.26'
$ws.Range("E27").Value = '  -0.41%  '
$ws.Range("E28").Value = '  -3.17%  '
$ws.Range("D29").Value = "'" + '18.96'
$ws.Range("E29").Value = '  -0.50%  '
$ws.Range("E30").Value = '  -6.06%  '
$ws.Range("E31").Value = '  +1.56%  '
$ws.Range("E32").Value = '  -2.19%  '
$ws.Range("E33").Value = '  +3.44%  '
$ws.Range("E34").Value = '  -2.17%  '
$ws.Range("D35").Value = "'" + '0.0603'
$ws.Range("E35").Value = '  -2.01%  '
$ws.Range("D36").Value = "'" + '6.46'
$ws.Range("E36").Value = '  +6.70%  '
$ws.Range("E37").Value = '  -5.18%  '
$ws.Range("E38").Value = '  -2.40%  '
$ws.Range("D39").Value = "'" + '1.00'
$ws.Range("E39").Value = '  +0.12%  '
$ws.Range("D40").Value = '1.540.57'
$ws.Range("E40").Value = '  +3.74%  '
$ws.Range("E41").Value = '  -0.12%  '
$ws.Range("B42").Value = 'Aave'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D42").Value = "'" + '97.19'
$ws.Range("E42").Value = '  -1.30%  '
$ws.Range("B43").Value = 'InjectiveProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D43").Value = "'" + '16.92'
$ws.Range("E43").Value = '  +0.40%  '
$ws.Range("E44").Value = '  -1.33%  '
$ws.Range("D45").Value = "'" + '0.0925'
$ws.Range("E45").Value = '  -2.07%  '
$ws.Range("E46").Value = '  -1.33%  '
$ws.Range("D47").Value = "'" + '3.93'
$ws.Range("E47").Value = '  -4.50%  '
$ws.Range("E48").Value = '  -1.26%  '
$ws.Range("D49").Value = "'" + '2.97'
$ws.Range("E49").Value = '  -0.35%  '
$ws.Range("E50").Value = '  +0.33%  '
$ws.Range("D51").Value = '2.231.56'
$ws.Range("E51").Value = '  -0.91%  '
